$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.258.39'
$ws.Cells.Item(2, 5).Value = '  +6.04%  '
$ws.Cells.Item(3, 4).Value = '2.246.08'
$ws.Cells.Item(3, 5).Value = '  +5.30%  '
$ws.Cells.Item(4, 5).Value = '  -0.17%  '
$ws.Cells.Item(5, 4).Value = "'252.63"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +8.20%  '
$ws.Cells.Item(6, 5).Value = '  +3.65%  '
$ws.Cells.Item(7, 4).Value = "'75.67"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +11.36%  '
$ws.Cells.Item(8, 5).Value = '  -0.24%  '
$ws.Cells.Item(9, 4).Value = "'0.605"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +9.23%  '
$ws.Cells.Item(10, 4).Value = "'41.29"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +9.99%  '
$ws.Cells.Item(11, 4).Value = "'0.0934"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +5.83%  '
$ws.Cells.Item(12, 4).Value = "'6.95"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +7.31%  '
$ws.Cells.Item(13, 5).Value = '  +3.05%  '
$ws.Cells.Item(14, 4).Value = '2.583.41'
$ws.Cells.Item(14, 5).Value = '  +5.25%  '
$ws.Cells.Item(15, 4).Value = "'14.66"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +3.33%  '
$ws.Cells.Item(16, 4).Value = '2.245.42'
$ws.Cells.Item(16, 5).Value = '  +5.74%  '
$ws.Cells.Item(17, 4).Value = "'0.794"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +3.96%  '
$ws.Cells.Item(18, 4).Value = '43.123.47'
$ws.Cells.Item(18, 5).Value = '  +6.10%  '
$ws.Cells.Item(19, 5).Value = '  +8.47%  '
$ws.Cells.Item(20, 4).Value = "'71.31"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +3.70%  '
$ws.Cells.Item(21, 4).Value = "'6.03"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +6.83%  '
$ws.Cells.Item(22, 4).Value = "'9.83"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +5.20%  '
$ws.Cells.Item(23, 5).Value = '  +20.38%  '
$ws.Cells.Item(24, 4).Value = "'230.23"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +4.06%  '
$ws.Cells.Item(25, 5).Value = '  +0.15%  '
$ws.Cells.Item(26, 4).Value = "'10.85"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +5.24%  '
$ws.Cells.Item(27, 5).Value = '  +7.09%  '
$ws.Cells.Item(28, 4).Value = "'2.28"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +8.34%  '
$ws.Cells.Item(29, 4).Value = "'39.34"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +32.40%  '
$ws.Cells.Item(30, 5).Value = '  +3.95%  '
$ws.Cells.Item(31, 4).Value = "'172.06"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +3.15%  '
$ws.Cells.Item(32, 4).Value = "'20.27"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +4.81%  '
$ws.Cells.Item(33, 4).Value = "'0.0804"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +9.42%  '
$ws.Cells.Item(34, 4).Value = "'5.30"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +7.10%  '
$ws.Cells.Item(35, 5).Value = '  +3.49%  '
$ws.Cells.Item(36, 5).Value = '  +13.39%  '
$ws.Cells.Item(37, 4).Value = "'4.48"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +13.23%  '
$ws.Cells.Item(38, 4).Value = "'0.0331"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +21.75%  '
$ws.Cells.Item(39, 4).Value = "'13.17"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +17.15%  '
$ws.Cells.Item(40, 5).Value = '  +6.09%  '
$ws.Cells.Item(41, 4).Value = "'0.206"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +14.19%  '
$ws.Cells.Item(42, 4).Value = "'5.44"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +4.80%  '
$ws.Cells.Item(43, 4).Value = "'59.72"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +6.37%  '
$ws.Cells.Item(44, 4).Value = "'104.90"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +10.80%  '
$ws.Cells.Item(45, 5).Value = '  +8.13%  '
$ws.Cells.Item(46, 4).Value = "'0.486"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +37.89%  '
$ws.Cells.Item(47, 4).Value = "'0.0997"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +6.39%  '
$ws.Cells.Item(48, 4).Value = "'2.41"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +15.15%  '
$ws.Cells.Item(49, 5).Value = '  +5.81%  '
$ws.Cells.Item(50, 4).Value = "'1.16"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +6.94%  '
$ws.Cells.Item(51, 4).Value = "'2.69"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +3.84%  '
